# Updated symbol list on Sun Dec 25 10:49:20 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (prices like "23.11" would
# otherwise be auto-coerced to a number by Excel), then restore the cell's
# style to the sheet's default "Normal" so no stray number-format sticks to it.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- simple price-only updates ---
Set-TextValue "D3" "23.11"
Set-TextValue "D5" "0.05977"
Set-TextValue "D8" "0.9282"

# --- rows 9-17: coin list rotated down by one; "One" moves from row 17 to row 9 ---
Set-TextValue "B9" "One"
Set-TextValue "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.01126"
Set-TextValue "E9" "8OneONEBestin24h"

Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1431"
Set-TextValue "E10" "9WazirXWRX"

Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07425"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"

Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03388"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"

Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03041"
Set-TextValue "E13" "12BitrueCoinBTR"

Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09338"
Set-TextValue "E14" "13BitMartTokenBMX"

Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.939"
Set-TextValue "E15" "14MCDexMCB"

Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001593"
Set-TextValue "E16" "15BitForexTokenBF"

Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04808"
Set-TextValue "E17" "16CoinExTokenCET"

# --- simple price-only updates ---
Set-TextValue "D18" "0.005545"
Set-TextValue "D19" "0.004156"
Set-TextValue "D20" "0.0009855"
Set-TextValue "D21" "0.00007705"
Set-TextValue "D22" "3.662"
Set-TextValue "D23" "6.459"
Set-TextValue "D26" "0.1349"
Set-TextValue "D27" "0.0002448"
Set-TextValue "D40" "0.03936"

# --- rows 41-43: coin list rotated down by one; "KickToken" moves from row 43 to row 41 ---
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006214"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1072"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002902"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D44" "0.007351"
Set-TextValue "E44" "43LocalTradersLCT"

Set-TextValue "D45" "0.00005136"
Set-TextValue "D47" "0.0005802"
Set-TextValue "D48" "0.8556"
Set-TextValue "D49" "0.002261"
